$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell and drop the stray "Exclude" type-label cell next to it.
$ws.Range("L1").Value = "#Exclude Table"
$ws.Range("M1").ClearContents()

# Reset the active selection to the top-left cell.
$ws.Range("L1").Select() | Out-Null
